$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1302.0849609375
$ws.Range("C2").Value = 0.9768
$ws.Range("D2").Value = 0.9413999915122986
$ws.Range("E2").Value = 1.496299982070923
$ws.Range("F2").Value = 0.7006000280380249
$ws.Range("H2").Value = 0.9457

# Row 3
$ws.Range("B3").Value = 1300.41357421875
$ws.Range("C3").Value = 1.0403
$ws.Range("D3").Value = 0.9539
$ws.Range("E3").Value = 2.096999883651733
$ws.Range("F3").Value = 0.6496999859809875
$ws.Range("H3").Value = 1.0569

# Row 4
$ws.Range("B4").Value = 864.810791015625
$ws.Range("C4").Value = 1.0103
$ws.Range("D4").Value = 0.9332
$ws.Range("E4").Value = 2.114099979400635
$ws.Range("F4").Value = 0.7580999732017517
$ws.Range("H4").Value = 0.8736

# Row 5
$ws.Range("B5").Value = 828.089599609375
$ws.Range("C5").Value = 0.8735000000000001
$ws.Range("D5").Value = 0.8643999999999999
$ws.Range("E5").Value = 1.305199980735779
$ws.Range("F5").Value = 0.6258000135421753
$ws.Range("H5").Value = 0.2639

# Row 6
$ws.Range("B6").Value = 1113.273193359375
$ws.Range("C6").Value = 0.8808
$ws.Range("D6").Value = 0.8718
$ws.Range("E6").Value = 1.135699987411499
$ws.Range("F6").Value = 0.6884999871253967
$ws.Range("H6").Value = 0.3296

# Row 7
$ws.Range("B7").Value = 853.4478759765625
$ws.Range("C7").Value = 0.8595
$ws.Range("D7").Value = 0.8547000288963318
$ws.Range("E7").Value = 1.060999989509583
$ws.Range("F7").Value = 0.7251999974250793
$ws.Range("H7").Value = 0.1779

# Row 8
$ws.Range("B8").Value = 919.458984375
$ws.Range("C8").Value = 0.8239
$ws.Range("D8").Value = 0.8184
$ws.Range("E8").Value = 1.119500041007996
$ws.Range("F8").Value = 0.7167999744415283
$ws.Range("H8").Value = -0.1433

# Row 9
$ws.Range("B9").Value = 7181.57861328125
$ws.Range("C9").Value = 0.9257
$ws.Range("D9").Value = 0.8848
$ws.Range("E9").Value = 2.114099979400635
$ws.Range("F9").Value = 0.6258000135421753
$ws.Range("H9").Value = 3.504300000000001
